# Fixing rules with LocalDate problem
#
# The "Set Queue Enter Date" rules (rows 27-31, column E) called
# toDate(java.time.LocalDate.now()) to convert the LocalDate into a
# java.util.Date before passing it to setQueueEnterDate. That conversion
# is unnecessary/incorrect here, so drop the toDate(...) wrapper and pass
# the LocalDate straight through.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E27:E31").Value = "setQueueEnterDate, java.time.LocalDate.now()"
